# Y4_B2526_General_&_Special_surgery_1_schedule.xlsx
# Change: every session's Duration (column G, rows 2-157) goes from 180 to 90,
# and the formatting of the "unshaded" rows in that column is consolidated to
# match the "shaded" rows (so the whole column ends up using a single style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 157
$col      = "G"

# Bulk-update all the duration values in one shot.
$ws.Range("$col$firstRow`:$col$lastRow").Value = 90

# The "shaded" style (gray fill) used by the even-numbered rows is the one we
# want every row in the column to end up with. Use it as the format source.
$shadedColor = 15790320
$src = $ws.Range("$col$firstRow")
$srcColor = $src.Interior.Color()
if ($srcColor -ne $shadedColor) {
    # fall back: find the first shaded row to copy from
    for ($r = $firstRow; $r -le $lastRow; $r++) {
        $cand = $ws.Range("$col$r")
        if ($cand.Interior.Color() -eq $shadedColor) {
            $src = $cand
            break
        }
    }
}
$src.Copy()

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Range("$col$r")
    if ($cell.Interior.Color() -ne $shadedColor) {
        $cell.PasteSpecial(-4122)
    }
}

$excel.CutCopyMode = $false

Write-Output "Updated $col$firstRow`:$col$lastRow to 90 and consolidated formatting."
